$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Colors_table")

# Row 4 corresponds to "Color Nr." = 3 (Magenta-Purple -> Magenta-Indigo)
$ws.Range("D4").Value = "#4B0082"
$ws.Range("B4").Value = "Magenta-Indigo"
$ws.Range("E4").Value = "https://github.com/Ing-Aladar-Dukay/CV_Dukay/blob/ca69b799c3560f1ac838964c7657905827b14bf9/03%20Colors%20icons/color%2003.png"
